# Working copy with excel and navigation changes
#
# - C2's hyperlink/email is changed from "ankita.sharma@gmail.com" to "Agasthya@gmail.com"
# - D2 gets a brand new hyperlink/email value "ankita.singh@gmail.com" (previously a
#   plain duplicate of C2's old value, without its own hyperlink)
# - The active selection on the sheet moves from J1 to C2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the existing hyperlink on C2 (its target is being replaced) -----------
$existingLinks = @($ws.Hyperlinks)
foreach ($link in $existingLinks) {
    if ($link.Range.Row -eq 2 -and $link.Range.Column -eq 3) {
        $link.Delete()
    }
}

# --- Update the cell values for C2 and D2 ------------------------------------------
$ws.Range("D2").Value = "ankita.singh@gmail.com"
$ws.Range("C2").Value = "Agasthya@gmail.com"

# --- (Re)create the hyperlinks: D2 is new, C2 points at the new address -----------
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:ankita.singh@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Agasthya@gmail.com") | Out-Null

# Make sure C2/D2 keep using the same shared "Hyperlink" cell style as E2
# (Hyperlinks.Add otherwise stamps its own duplicate style record on the cell)
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("D2").Style = "Hyperlink"

# --- Update the selected / active cell on the sheet --------------------------------
$ws.Range("C2").Select() | Out-Null
